$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "Datos actualizados..." timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 22:59"

# Update Cataluña row (row 5) figures
$ws.Range("B5").Value = 15026
$ws.Range("C5").Value = 3455
$ws.Range("D5").Value = 10345
$ws.Range("E5").Value = 1226
